$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 5127.5
$ws.Range("I11").Value = 5127.5
$ws.Range("K11").Value = 5127.5
$ws.Range("M11").Value = -4987.5

$ws.Range("H99").Value = 403.76923
$ws.Range("I99").Value = 379.16666
$ws.Range("J99").Value = 699
$ws.Range("K99").Value = 1137.49998
$ws.Range("L99").Value = 2097
$ws.Range("M99").Value = 360.5000199999999
$ws.Range("N99").Value = -5093

$ws.Range("H132").Value = 1575.4412
$ws.Range("I132").Value = 966.3929000000001
$ws.Range("J132").Value = 4417.6665
$ws.Range("K132").Value = 2899.1787
$ws.Range("L132").Value = 13252.9995
$ws.Range("M132").Value = -369.1787000000004
$ws.Range("N132").Value = -18312.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 978.45
$ws.Range("I2").Value = 870.4666999999999
$ws.Range("K2").Value = 870.4666999999999
$ws.Range("M2").Value = -757.4666999999999

$ws.Range("H31").Value = 25987
$ws.Range("J31").Value = 111500
$ws.Range("L31").Value = 111500
$ws.Range("N31").Value = -112088

$ws.Range("H32").Value = 16136786
$ws.Range("I32").Value = 21741756
$ws.Range("J32").Value = 22501.375
$ws.Range("K32").Value = 21741756
$ws.Range("L32").Value = 22501.375
$ws.Range("M32").Value = -21741469
$ws.Range("N32").Value = -23075.375

$ws.Range("H45").Value = 20836318
$ws.Range("I45").Value = 35716300
$ws.Range("K45").Value = 35716300
$ws.Range("M45").Value = -35715923

$ws.Range("H61").Value = 41674310
$ws.Range("I61").Value = 83339980
$ws.Range("K61").Value = 83339980
$ws.Range("M61").Value = -83339768

$ws.Range("H110").Value = 4333.3335
$ws.Range("I110").Value = 4331
$ws.Range("J110").Value = 4335.6665
$ws.Range("K110").Value = 4331
$ws.Range("L110").Value = 4335.6665
$ws.Range("M110").Value = -2286
$ws.Range("N110").Value = -8425.666499999999

$ws.Range("H116").Value = 978.45
$ws.Range("I116").Value = 870.4666999999999
$ws.Range("K116").Value = 870.4666999999999
$ws.Range("M116").Value = 1423.5333

$ws.Range("H136").Value = 41674310
$ws.Range("I136").Value = 83339980
$ws.Range("K136").Value = 250019940
$ws.Range("M136").Value = -250017390

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 978.45
$ws.Range("I3").Value = 870.4666999999999
$ws.Range("K3").Value = 870.4666999999999
$ws.Range("M3").Value = -756.4666999999999

$ws.Range("H105").Value = 4676.826
$ws.Range("I105").Value = 13534
$ws.Range("K105").Value = 13534
$ws.Range("M105").Value = -11787

$ws.Range("H107").Value = 2949.8667
$ws.Range("J107").Value = 2487.4
$ws.Range("L107").Value = 2487.4
$ws.Range("N107").Value = -6327.4

$ws.Range("H134").Value = 426428.94
$ws.Range("I134").Value = 1714.875
$ws.Range("J134").Value = 1105971.4
$ws.Range("K134").Value = 5144.625
$ws.Range("L134").Value = 3317914.2
$ws.Range("M134").Value = -2609.625
$ws.Range("N134").Value = -3322984.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1366.8334

$ws.Range("H109").Value = 39999.5
$ws.Range("J109").Value = 39999.5
$ws.Range("L109").Value = 39999.5
$ws.Range("N109").Value = -42079.5

$ws.Range("H113").Value = 1366.8334

$ws.Range("H132").Value = 4153.5557
$ws.Range("I132").Value = 2652.1333
$ws.Range("K132").Value = 7956.3999
$ws.Range("M132").Value = -5426.3999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 162.55882
$ws.Range("J2").Value = 201.86363
$ws.Range("L2").Value = 1211.18178
$ws.Range("N2").Value = -1437.18178

$ws.Range("H37").Value = 58624.25
$ws.Range("J37").Value = 58624.25
$ws.Range("L37").Value = 175872.75
$ws.Range("N37").Value = -176096.75

$ws.Range("H131").Value = 9734.026
$ws.Range("J131").Value = 9734.026
$ws.Range("L131").Value = 29202.078
$ws.Range("N131").Value = -39282.078

$ws.Range("H134").Value = 9846.517
$ws.Range("J134").Value = 12657.087
$ws.Range("L134").Value = 37971.261
$ws.Range("N134").Value = -48111.261

$ws.Range("H137").Value = 4839.615
$ws.Range("I137").Value = 4120
$ws.Range("K137").Value = 12360
$ws.Range("M137").Value = -7260

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 680.6667
$ws.Range("I107").Value = 522
$ws.Range("J107").Value = 998
$ws.Range("K107").Value = 522
$ws.Range("L107").Value = 998
$ws.Range("M107").Value = 1398
$ws.Range("N107").Value = -4838

$ws.Range("H132").Value = 52641744
$ws.Range("I132").Value = 83337560
$ws.Range("J132").Value = 20344.285
$ws.Range("K132").Value = 250012680
$ws.Range("L132").Value = 61032.855
$ws.Range("M132").Value = -250010150
$ws.Range("N132").Value = -66092.855

$ws.Range("H134").Value = 101449.75
$ws.Range("J134").Value = 101449.75
$ws.Range("L134").Value = 304349.25
$ws.Range("N134").Value = -309419.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1316.8889
$ws.Range("I22").Value = 1341.6666
$ws.Range("J22").Value = 1267.3334
$ws.Range("K22").Value = 1341.6666
$ws.Range("L22").Value = 1267.3334
$ws.Range("M22").Value = -1046.6666
$ws.Range("N22").Value = -1857.3334

$ws.Range("H27").Value = 1316.8889
$ws.Range("I27").Value = 1341.6666
$ws.Range("J27").Value = 1267.3334
$ws.Range("K27").Value = 1341.6666
$ws.Range("L27").Value = 1267.3334
$ws.Range("M27").Value = -1234.6666
$ws.Range("N27").Value = -1481.3334

$ws.Range("H68").Value = 2073
$ws.Range("I68").Value = 1999
$ws.Range("J68").Value = 2110
$ws.Range("K68").Value = 1999
$ws.Range("L68").Value = 2110
$ws.Range("M68").Value = -1250
$ws.Range("N68").Value = -3608

$ws.Range("H71").Value = 2073
$ws.Range("I71").Value = 1999
$ws.Range("J71").Value = 2110
$ws.Range("K71").Value = 9995
$ws.Range("L71").Value = 10550
$ws.Range("M71").Value = -6251
$ws.Range("N71").Value = -18038

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 6016
$ws.Range("I81").Value = 1910.3334
$ws.Range("K81").Value = 3820.6668
$ws.Range("M81").Value = -2759.6668

$ws.Range("H84").Value = 6016
$ws.Range("I84").Value = 1910.3334
$ws.Range("K84").Value = 19103.334
$ws.Range("M84").Value = -13799.334

$ws.Range("H93").Value = 109000
$ws.Range("J93").Value = 109000
$ws.Range("L93").Value = 109000
$ws.Range("N93").Value = -113992

$ws.Range("H107").Value = 31250862
$ws.Range("I107").Value = 38462336
$ws.Range("K107").Value = 115387008
$ws.Range("M107").Value = -115385088

$ws.Range("H132").Value = 306299.56
$ws.Range("J132").Value = 2003741.4
$ws.Range("L132").Value = 6011224.199999999
$ws.Range("N132").Value = -6016284.199999999
